# Apply cryptos list update (GitHub Actions refresh) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, [string]$CellRef, [string]$NewValue)
    $c = $Worksheet.Range($CellRef)
    # Force the cell to keep its literal text representation (avoid Excel
    # auto-converting number-like strings such as "307.98" into numerics).
    $c.NumberFormat = "@"
    $c.Value = $NewValue
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '43.034.71'
$ws.Range('E2').Value = '  -0.16%  '
Set-TextCell $ws 'D3' '2.354.12'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  +0.15%  '
Set-TextCell $ws 'D5' '307.98'
$ws.Range('E5').Value = '  -0.66%  '
Set-TextCell $ws 'D6' '102.41'
$ws.Range('E6').Value = '  +1.29%  '
Set-TextCell $ws 'D7' '0.509'
$ws.Range('E7').Value = '  -4.84%  '
$ws.Range('E8').Value = '  +0.12%  '
Set-TextCell $ws 'D9' '0.516'
$ws.Range('E9').Value = '  -1.59%  '
Set-TextCell $ws 'D10' '35.29'
$ws.Range('E10').Value = '  -2.03%  '
Set-TextCell $ws 'D11' '52.57'
$ws.Range('E11').Value = '  +0.81%  '
Set-TextCell $ws 'D12' '0.0801'
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('E13').Value = '  -0.56%  '
Set-TextCell $ws 'D14' '6.85'
$ws.Range('E14').Value = '  -3.82%  '
Set-TextCell $ws 'D15' '2.731.16'
$ws.Range('E15').Value = '  +2.20%  '
Set-TextCell $ws 'D16' '15.44'
$ws.Range('E16').Value = '  +3.06%  '
Set-TextCell $ws 'D17' '2.362.74'
$ws.Range('E17').Value = '  +2.13%  '
Set-TextCell $ws 'D18' '0.803'
$ws.Range('E18').Value = '  -1.32%  '
Set-TextCell $ws 'D19' '43.067.45'
$ws.Range('E19').Value = '  +0.09%  '
Set-TextCell $ws 'D20' '6.26'
$ws.Range('E20').Value = '  +2.24%  '
Set-TextCell $ws 'D21' '11.74'
$ws.Range('E21').Value = '  -6.49%  '
Set-TextCell $ws 'D22' '0.0₃0908'
$ws.Range('E22').Value = '  -1.31%  '
Set-TextCell $ws 'D23' '67.77'
$ws.Range('E23').Value = '  -0.93%  '
Set-TextCell $ws 'D24' '238.14'
$ws.Range('E24').Value = '  -1.22%  '
Set-TextCell $ws 'D25' '2.01'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('E27').Value = '  +0.07%  '
Set-TextCell $ws 'D28' '25.42'
$ws.Range('E28').Value = '  +2.97%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws 'D29' '3.86'
$ws.Range('E29').Value = '  -3.03%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D30' '2.20'
$ws.Range('E30').Value = '  +4.11%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D31' '35.87'
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D32' '9.37'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D33' '160.89'
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D35' '5.17'
$ws.Range('E35').Value = '  -3.06%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell $ws 'D36' '17.86'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws 'D37' '2.49'
$ws.Range('E37').Value = '  +4.23%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws 'D38' '4.64'
$ws.Range('E38').Value = '  +8.35%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D39' '3.04'
$ws.Range('E39').Value = '  -4.07%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D40' '0.0729'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D41' '1.89'
$ws.Range('E41').Value = '  +2.52%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D42' '0.104'
$ws.Range('E42').Value = '  -3.12%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D43' '0.113'
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws 'D44' '2.59'
$ws.Range('E44').Value = '  +11.95%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D45' '2.029.85'
$ws.Range('E45').Value = '  +2.81%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D46' '19.44'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D47' '0.0287'
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws 'D48' '10.54'
$ws.Range('E48').Value = '  +7.48%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D49' '3.03'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell $ws 'D50' '57.18'
$ws.Range('E50').Value = '  +2.54%  '
Set-TextCell $ws 'D51' '2.593.46'
$ws.Range('E51').Value = '  +2.12%  '
